$wb = $excel.ActiveWorkbook

# Sheet "展览" (1st sheet): update F2, F3, F4
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1399
$ws1.Range("F3").Value = 2944
$ws1.Range("F4").Value = 20

# Sheet "全部类型" (4th sheet): update F3, F4, F5
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1399
$ws4.Range("F4").Value = 2944
$ws4.Range("F5").Value = 20
